# daily auto push: 2026-02-21 04:12 UTC
#
# A new sample row (2026/02/21, 土, 12, 201) was recorded between the
# existing "2026/02/21" rows (row 851) and the old "2026/12/29" row
# (old row 852). Inserting a whole row at sheet row 852 shifts every
# subsequent row down by one (old 852-893 -> new 853-894), which matches
# the rest of the diff exactly (each later row's data is simply the
# previous row's data, one row lower), and grows the used range from
# A1:D893 to A1:D894.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push everything at/after row 852 down by one row.
$ws.Rows(852).Insert()

# Populate the newly-opened row 852 with the new sample.
# Copy the date/weekday text from the row above (itself "2026/02/21" / "土")
# so the new cells stay plain text instead of being reinterpreted as a
# date serial when assigned as a string literal.
$ws.Range("A851:B851").Copy($ws.Range("A852:B852"))
$ws.Range("C852").Value = 12
$ws.Range("D852").Value = 201
